$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H column labels to pluralized/expanded subcategories
$ws.Range("H3").Value = "line graph(s)"
$ws.Range("H5").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H9").Value = "line graph(s)"
$ws.Range("H10").Value = "line graph(s)"
$ws.Range("H11").Value = "line graph(s)"
$ws.Range("H12").Value = "line graph(s)"
$ws.Range("H13").Value = "line graph(s)"
$ws.Range("H14").Value = "line graph(s)"
$ws.Range("H15").Value = "bar chart(s)"
$ws.Range("H21").Value = "line graph(s)"

# Remove the "is_viewed" column entirely (column I)
$ws.Columns.Item(9).Delete()
